# Attendance.xlsx edit: change H2 from a text percentage ("47.62%") to the
# plain numeric value 47.62, then add a new blank worksheet ("Sheet1") and
# make it the active/selected tab, with the cursor on Attendance moved to H9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

# H2 used to be the text "47.62%" stored via shared string; make it a real number.
$ws.Range("H2").Value = 47.62

# Give column H an explicit width like the other sheets in the diff (closest
# reachable width to the authored 11.08203125).
$ws.Columns.Item(8).ColumnWidth = 10.1

# Move the selection on the Attendance sheet from I1 to H9 (matches the diff).
$ws.Range("H9").Select()

# Add a new (blank) worksheet named "Sheet1" right after the Attendance sheet;
# it becomes the active sheet/tab.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Sheet1"
